# "colors more standardized within ptfiltered mod"
#
# AntennaMetadata: fill in the new PressureTransducerSiteName column (E)
# for the Connectivity Channel stationary-antenna rows, and move the
# sheet selection to that column's header cell.
$wb = $excel.ActiveWorkbook

$antenna = $wb.Worksheets.Item("AntennaMetadata")
$antenna.Range("E13").Value = "Connectivity Downstream"
$antenna.Range("E14").Value = "Connectivity Downstream"
$antenna.Range("E15").Value = "Connectivity Side Channel"
$antenna.Range("E16").Value = "Connectivity Side Channel"
$antenna.Range("E17").Value = "Connectivity Upstream"
$antenna.Range("E18").Value = "Connectivity Upstream"
$antenna.Activate()
$antenna.Range("E1").Select()

# Notes: add a reminder about keeping PressureTransducerSiteName lined up
# with the pressure-transducer data, and leave that sheet active/selected.
$notes = $wb.Worksheets.Item("Notes")
$notes.Range("A6").Value = "PressureTransducerSiteName needs to line up with the sites that are in the Pressuretransducer data"
$notes.Activate()
$notes.Range("M5").Select()
